$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17,8).Value = 3333.1667
$ws.Cells.Item(17,10).Value = 3333.1667
$ws.Cells.Item(17,12).Value = 9999.500100000001
$ws.Cells.Item(17,14).Value = -10335.5001
$ws.Cells.Item(18,8).Value = 7025
$ws.Cells.Item(18,9).Value = 7025
$ws.Cells.Item(18,10).Value = 0
$ws.Cells.Item(18,11).Value = 7025
$ws.Cells.Item(18,12).Value = 0
$ws.Cells.Item(18,13).Value = -6741
$ws.Cells.Item(18,14).ClearContents()
$ws.Cells.Item(111,8).Value = 989
$ws.Cells.Item(111,9).Value = 989
$ws.Cells.Item(111,10).Value = 0
$ws.Cells.Item(111,11).Value = 2967
$ws.Cells.Item(111,12).Value = 0
$ws.Cells.Item(111,13).Value = 100
$ws.Cells.Item(111,14).ClearContents()
$ws.Cells.Item(112,8).Value = 969.2
$ws.Cells.Item(112,10).Value = 969.2
$ws.Cells.Item(112,12).Value = 2907.6
$ws.Cells.Item(112,14).Value = -5123.6
$ws.Cells.Item(138,8).Value = 3835.0342
$ws.Cells.Item(138,10).Value = 4104.7075
$ws.Cells.Item(138,12).Value = 12314.1225
$ws.Cells.Item(138,14).Value = -22594.1225
$ws.Cells.Item(141,8).Value = 4780.1665
$ws.Cells.Item(141,9).Value = 4346.3
$ws.Cells.Item(141,11).Value = 13038.9
$ws.Cells.Item(141,13).Value = -7858.900000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(14,8).Value = 0
$ws.Cells.Item(14,9).Value = 0
$ws.Cells.Item(14,10).Value = 0
$ws.Cells.Item(14,11).Value = 0
$ws.Cells.Item(14,12).Value = 0
$ws.Cells.Item(14,13).ClearContents()
$ws.Cells.Item(14,14).ClearContents()
$ws.Cells.Item(32,8).Value = 4975.357
$ws.Cells.Item(32,9).Value = 4435.0386
$ws.Cells.Item(32,11).Value = 4435.0386
$ws.Cells.Item(32,13).Value = -4148.0386
$ws.Cells.Item(74,8).Value = 2306.7144
$ws.Cells.Item(74,9).Value = 1935.6666
$ws.Cells.Item(74,11).Value = 1935.6666
$ws.Cells.Item(74,13).Value = -1061.6666
$ws.Cells.Item(77,8).Value = 2306.7144
$ws.Cells.Item(77,9).Value = 1935.6666
$ws.Cells.Item(77,11).Value = 9678.333000000001
$ws.Cells.Item(77,13).Value = -5310.333000000001
$ws.Cells.Item(88,8).Value = 1645.5834
$ws.Cells.Item(88,10).Value = 1769.7
$ws.Cells.Item(88,12).Value = 1769.7
$ws.Cells.Item(88,14).Value = -2581.7
$ws.Cells.Item(91,8).Value = 1645.5834
$ws.Cells.Item(91,10).Value = 1769.7
$ws.Cells.Item(91,12).Value = 1769.7
$ws.Cells.Item(91,14).Value = -4577.7
$ws.Cells.Item(102,8).Value = 1331.4445
$ws.Cells.Item(102,9).Value = 1395.0588
$ws.Cells.Item(102,11).Value = 1395.0588
$ws.Cells.Item(102,13).Value = 226.9412
$ws.Cells.Item(110,8).Value = 1194.2222
$ws.Cells.Item(110,9).Value = 1194.2222
$ws.Cells.Item(110,11).Value = 1194.2222
$ws.Cells.Item(110,13).Value = 850.7778000000001
$ws.Cells.Item(132,8).Value = 3058.2
$ws.Cells.Item(132,9).Value = 3058.2
$ws.Cells.Item(132,11).Value = 9174.599999999999
$ws.Cells.Item(132,13).Value = -6644.599999999999
$ws.Cells.Item(139,8).Value = 0
$ws.Cells.Item(139,10).Value = 0
$ws.Cells.Item(139,12).Value = 0
$ws.Cells.Item(139,14).ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(96,8).Value = 8868.444
$ws.Cells.Item(96,9).Value = 8868.444
$ws.Cells.Item(96,11).Value = 8868.444
$ws.Cells.Item(96,13).Value = -6122.444
$ws.Cells.Item(107,8).Value = 3957.5
$ws.Cells.Item(107,9).Value = 3822.5789
$ws.Cells.Item(107,11).Value = 3822.5789
$ws.Cells.Item(107,13).Value = -1902.5789
$ws.Cells.Item(133,8).Value = 69999
$ws.Cells.Item(133,9).Value = 0
$ws.Cells.Item(133,11).Value = 0
$ws.Cells.Item(133,13).ClearContents()
$ws.Cells.Item(134,8).Value = 977.6667
$ws.Cells.Item(134,9).Value = 975.6
$ws.Cells.Item(134,10).Value = 988
$ws.Cells.Item(134,11).Value = 2926.8
$ws.Cells.Item(134,12).Value = 2964
$ws.Cells.Item(134,13).Value = -391.8000000000002
$ws.Cells.Item(134,14).Value = -8034

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31,8).Value = 7902.8335
$ws.Cells.Item(31,9).Value = 15997
$ws.Cells.Item(31,10).Value = 7167
$ws.Cells.Item(31,11).Value = 15997
$ws.Cells.Item(31,12).Value = 7167
$ws.Cells.Item(31,13).Value = -15702
$ws.Cells.Item(31,14).Value = -7757
$ws.Cells.Item(34,8).Value = 7902.8335
$ws.Cells.Item(34,9).Value = 15997
$ws.Cells.Item(34,10).Value = 7167
$ws.Cells.Item(34,11).Value = 15997
$ws.Cells.Item(34,12).Value = 7167
$ws.Cells.Item(34,13).Value = -15795
$ws.Cells.Item(34,14).Value = -7571
$ws.Cells.Item(45,8).Value = 67
$ws.Cells.Item(45,9).Value = 67
$ws.Cells.Item(45,11).Value = 67
$ws.Cells.Item(45,13).Value = 526
$ws.Cells.Item(58,8).Value = 881.7143
$ws.Cells.Item(58,9).Value = 881.7143
$ws.Cells.Item(58,11).Value = 881.7143
$ws.Cells.Item(58,13).Value = -678.7143
$ws.Cells.Item(132,8).Value = 1901.2354
$ws.Cells.Item(132,9).Value = 1540.8
$ws.Cells.Item(132,11).Value = 4622.4
$ws.Cells.Item(132,13).Value = -2092.4
$ws.Cells.Item(134,8).Value = 2761.875
$ws.Cells.Item(134,9).Value = 2783
$ws.Cells.Item(134,11).Value = 8349
$ws.Cells.Item(134,13).Value = -5814
$ws.Cells.Item(136,8).Value = 881.7143
$ws.Cells.Item(136,9).Value = 881.7143
$ws.Cells.Item(136,11).Value = 2645.1429
$ws.Cells.Item(136,13).Value = -95.14289999999983

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5,8).Value = 4576.857
$ws.Cells.Item(5,9).Value = 5292.1665
$ws.Cells.Item(5,10).Value = 285
$ws.Cells.Item(5,11).Value = 15876.4995
$ws.Cells.Item(5,12).Value = 855
$ws.Cells.Item(5,13).Value = -15764.4995
$ws.Cells.Item(5,14).Value = -1079
$ws.Cells.Item(75,8).Value = 3066.182
$ws.Cells.Item(75,10).Value = 3341.2856
$ws.Cells.Item(75,12).Value = 10023.8568
$ws.Cells.Item(75,14).Value = -12019.8568
$ws.Cells.Item(78,8).Value = 3066.182
$ws.Cells.Item(78,10).Value = 3341.2856
$ws.Cells.Item(78,12).Value = 30071.5704
$ws.Cells.Item(78,14).Value = -40055.5704
$ws.Cells.Item(115,8).Value = 400
$ws.Cells.Item(115,9).Value = 400
$ws.Cells.Item(115,11).Value = 1200
$ws.Cells.Item(115,13).Value = -25
$ws.Cells.Item(131,8).Value = 2799.6155
$ws.Cells.Item(131,9).Value = 1998.75
$ws.Cells.Item(131,10).Value = 3155.5557
$ws.Cells.Item(131,11).Value = 5996.25
$ws.Cells.Item(131,12).Value = 9466.667099999999
$ws.Cells.Item(131,13).Value = -956.25
$ws.Cells.Item(131,14).Value = -19546.6671
$ws.Cells.Item(135,8).Value = 4576.857
$ws.Cells.Item(135,9).Value = 5292.1665
$ws.Cells.Item(135,10).Value = 285
$ws.Cells.Item(135,11).Value = 47629.4985
$ws.Cells.Item(135,12).Value = 2565
$ws.Cells.Item(135,13).Value = -45094.4985
$ws.Cells.Item(135,14).Value = -7635

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(18,8).Value = 4154.8335
$ws.Cells.Item(107,8).Value = 318.42856
$ws.Cells.Item(107,9).Value = 271.5
$ws.Cells.Item(107,11).Value = 271.5
$ws.Cells.Item(107,13).Value = 1648.5
$ws.Cells.Item(126,8).Value = 4998
$ws.Cells.Item(126,9).Value = 4872
$ws.Cells.Item(126,10).Value = 5250
$ws.Cells.Item(126,11).Value = 14616
$ws.Cells.Item(126,12).Value = 15750
$ws.Cells.Item(126,13).Value = -12146
$ws.Cells.Item(126,14).Value = -20690
$ws.Cells.Item(132,8).Value = 1744
$ws.Cells.Item(132,9).Value = 1548.6471
$ws.Cells.Item(132,11).Value = 4645.9413
$ws.Cells.Item(132,13).Value = -2115.9413

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(12,8).Value = 0
$ws.Cells.Item(12,10).Value = 0
$ws.Cells.Item(12,12).Value = 0
$ws.Cells.Item(12,14).ClearContents()
$ws.Cells.Item(22,8).Value = 2861.2
$ws.Cells.Item(22,9).Value = 3468.75
$ws.Cells.Item(22,11).Value = 3468.75
$ws.Cells.Item(22,13).Value = -3173.75
$ws.Cells.Item(27,8).Value = 2861.2
$ws.Cells.Item(27,9).Value = 3468.75
$ws.Cells.Item(27,11).Value = 3468.75
$ws.Cells.Item(27,13).Value = -3361.75
$ws.Cells.Item(61,8).Value = 5433
$ws.Cells.Item(61,10).Value = 6952
$ws.Cells.Item(61,12).Value = 6952
$ws.Cells.Item(61,14).Value = -7356
$ws.Cells.Item(113,8).Value = 5433
$ws.Cells.Item(113,10).Value = 6952
$ws.Cells.Item(113,12).Value = 6952
$ws.Cells.Item(113,14).Value = -11292
$ws.Cells.Item(125,8).Value = 0
$ws.Cells.Item(125,10).Value = 0
$ws.Cells.Item(125,12).Value = 0
$ws.Cells.Item(125,14).ClearContents()
$ws.Cells.Item(141,8).Value = 41101.832
$ws.Cells.Item(141,9).Value = 22000
$ws.Cells.Item(141,10).Value = 79305.5
$ws.Cells.Item(141,11).Value = 22000
$ws.Cells.Item(141,12).Value = 79305.5
$ws.Cells.Item(141,13).Value = -16820
$ws.Cells.Item(141,14).Value = -89665.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(99,8).Value = 0
$ws.Cells.Item(99,9).Value = 0
$ws.Cells.Item(99,11).Value = 0
$ws.Cells.Item(99,13).ClearContents()
$ws.Cells.Item(132,8).Value = 3337.8572
$ws.Cells.Item(132,9).Value = 3337.8572
$ws.Cells.Item(132,11).Value = 10013.5716
$ws.Cells.Item(132,13).Value = -7483.571599999999
